$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.793.88'
$ws.Range("E2").Value = '  -2.42%  '
$ws.Range("D3").Value = '1.569.04'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("E6").Value = '  -2.23%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.02'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("E10").Value = '  -1.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("D12").Value = '1.792.23'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '1.575.53'
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("E14").Value = '  -2.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.515'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("D16").Value = '26.800.43'
$ws.Range("E16").Value = '  -2.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.44'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.62%  '
$ws.Range("E18").Value = '  +1.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.32'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").Value = '0.0₃0677'
$ws.Range("E20").Value = '  -1.91%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.49%  '
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("E29").Value = '  -1.38%  '
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("E31").Value = '  -3.34%  '
$ws.Range("E32").Value = '  -1.30%  '
$ws.Range("D33").Value = '1.396.06'
$ws.Range("E33").Value = '  +1.23%  '
$ws.Range("E34").Value = '  -1.26%  '
$ws.Range("E35").Value = '  -0.96%  '
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.934'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.98%  '
$ws.Range("E38").Value = '  -2.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.530'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.817'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.43%  '
$ws.Range("D47").Value = '1.704.92'
$ws.Range("E47").Value = '  +0.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.82%  '
$ws.Range("D49").Value = '0.0₇0986'
$ws.Range("E49").Value = '  -1.35%  '
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("E51").Value = '  -0.74%  '
